# Supplementary tables workbook update:
#  - rename the two existing sheets to include their new "Table S#" labels
#  - add captions (bold "Table S#" lead-in + description) under each table
#  - add a brand new "Table S3, Better or worse" sheet with the polisher
#    better/same/worse counts and its own caption

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename existing sheets
# ---------------------------------------------------------------------------
$wsAccessions = $wb.Worksheets.Item("Accessions")
$wsAccessions.Name = "Table S1, accessions"

$wsMedaka = $wb.Worksheets.Item("Medaka polishing")
$wsMedaka.Name = "Table S2, Medaka polishing"

# ---------------------------------------------------------------------------
# 2) Add caption under Table S1 (accessions)
# ---------------------------------------------------------------------------
$s1CaptionRow = 12
$wsAccessions.Range("A$s1CaptionRow`:D$s1CaptionRow").Merge() | Out-Null
$s1Cell = $wsAccessions.Range("A$s1CaptionRow")
$s1Text = "Table S1: NCBI accessions for each of the nine genomes used in this study and their ONT and Illumina read sets."
$s1Cell.Value = $s1Text
$s1Cell.Characters(1, 8).Font.Bold = $true
$wsAccessions.Range("A$s1CaptionRow`:D$s1CaptionRow").WrapText = $true
$wsAccessions.Range("A$s1CaptionRow`:D$s1CaptionRow").VerticalAlignment = -4108
$wsAccessions.Rows.Item($s1CaptionRow).RowHeight = 44

# ---------------------------------------------------------------------------
# 3) Update the caption under Table S2 (Medaka polishing) with the new
#    "Table S2:" lead-in and the "(blue)"/"(red)" clarifications
# ---------------------------------------------------------------------------
$s2CaptionRow = 14
$s2Cell = $wsMedaka.Range("A$s2CaptionRow")
$s2Text = "Table S2: Of the nine ONT-only assemblies used in this study, three improved with Medaka polishing (blue), three did not change, three got worse (red), and Medaka made the total error count worse. We therefore do not recommend using Medaka to polish Trycycler assemblies of sup-basecalled ONT reads."
$s2Cell.Value = $s2Text
$s2Cell.Characters(1, 8).Font.Bold = $true
$wsMedaka.Rows.Item($s2CaptionRow).RowHeight = 71

# ---------------------------------------------------------------------------
# 4) Add the new "Table S3, Better or worse" sheet at the end of the workbook
# ---------------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws3 = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws3.Name = "Table S3, Better or worse"

# Header row
$ws3.Range("A1").Value = "Polisher"
$ws3.Range("B1").Value = "Better"
$ws3.Range("C1").Value = "Same"
$ws3.Range("D1").Value = "Worse"

# Copy the header formatting (bold, centered, bottom border) from Table S1's
# header row so the new table matches the rest of the workbook.
$wsAccessions.Range("A1:D1").Copy() | Out-Null
$ws3.Range("A1:D1").PasteSpecial(-4122) | Out-Null
$ws3.Rows.Item(1).RowHeight = 18

# Data rows: polisher name, better count, same count, worse count
$polisherData = @(
    @("Polypolish default", 3260, 1093, 147),
    @("Polypolish careful", 3321, 1179, 0),
    @("Pypolca default", 2656, 470, 1374),
    @("Pypolca careful", 3526, 923, 51),
    @("HyPo", 1007, 396, 3097),
    @("FMLRC2", 1450, 379, 2671),
    @("NextPolish", 1851, 460, 2189),
    @("Pilon", 2462, 707, 1331)
)

$row = 2
foreach ($rec in $polisherData) {
    $ws3.Cells.Item($row, 1).Value = $rec[0]
    $ws3.Cells.Item($row, 2).Value = $rec[1]
    $ws3.Cells.Item($row, 3).Value = $rec[2]
    $ws3.Cells.Item($row, 4).Value = $rec[3]
    $row++
}

# Last data row (Pilon, row 9) gets the closing bottom border, copied from
# Table S1's final row so the style matches exactly.
$lastDataRow = $row - 1
$wsAccessions.Range("A10:D10").Copy() | Out-Null
$ws3.Range("A$lastDataRow`:D$lastDataRow").PasteSpecial(-4122) | Out-Null

# Column widths
$ws3.Columns.Item(1).ColumnWidth = 14.830729166666666
$ws3.Range("B1:D1").ColumnWidth = 8.666666666666666

# Caption row
$s3CaptionRow = 11
$ws3.Range("A$s3CaptionRow`:D$s3CaptionRow").Merge() | Out-Null
$s3Cell = $ws3.Range("A$s3CaptionRow")
$s3Text = "Table S3: Each polisher was run 4500 times (9 genomes at 500 depths). Each result was classified as 'Better' (fewer errors after polishing), 'Same' (no change in errors after polishing) or 'Worse' (more errors after polishing), with the totals shown in this table."
$s3Cell.Value = $s3Text
$s3Cell.Characters(1, 8).Font.Bold = $true
$ws3.Range("A$s3CaptionRow`:D$s3CaptionRow").WrapText = $true
$ws3.Range("A$s3CaptionRow`:D$s3CaptionRow").VerticalAlignment = -4108
$ws3.Rows.Item($s3CaptionRow).RowHeight = 104

# Leave the new sheet's selection on A1
$ws3.Range("A1").Select() | Out-Null
